$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix a handful of mistyped Note_GINF31 (column D) grades.
$ws.Range("D4").Value = 12
$ws.Range("D5").Value = 13
$ws.Range("D6").Value = 14
$ws.Range("D8").Value = 9
$ws.Range("D9").Value = 10
$ws.Range("D10").Value = 7
$ws.Range("D12").Value = 11
$ws.Range("D14").Value = 17

# Remove the MOY (average) column entirely - header + per-row AVERAGE formulas.
$ws.Columns("P").Delete()

# Update the view: scroll so column E is the first visible column, and select P1
# (the first empty column after the table, now that the MOY column is gone).
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("P1").Select()
